$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E1 (exponent) from 1.6 to 1
$ws.Range("E1").Value = 1

# Update D1 formula to reference C1 instead of SIN(B1)
$ws.Range("D1").Formula = '=POWER(C1,$E$1)'

# Update the shared-formula master in D2; D3:D51 inherit the new formula automatically
$ws.Range("D2").Formula = '=POWER(C2,$E$1)'

# Move selection to I29 like in the diff
$ws.Range("I29").Select()
